$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested-count) values on both the "展览" sheet
# and the consolidated "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5193
    $ws.Range("F3").Value = 159
    $ws.Range("F4").Value = 909
}
